# This script corrects a set of rows (52-57, 59-60) in the active worksheet
# whose data got shuffled: each affected row currently holds the record that
# belongs to a different row. We restore the correct record per row by
# permuting columns A, B, D, E, F, G, H, P, Q, R (identity/taxon columns)
# across the affected rows; all other columns already match and are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot current values (must read everything before any writes, since
# the target rows form permutation cycles and would otherwise clobber source
# data before it is read). ---
$A52 = $ws.Range("A52").Value2
$B52 = $ws.Range("B52").Value2
$D52 = $ws.Range("D52").Value2
$E52 = $ws.Range("E52").Value2
$F52 = $ws.Range("F52").Value2
$G52 = $ws.Range("G52").Value2
$H52 = $ws.Range("H52").Value2
$P52 = $ws.Range("P52").Value2
$Q52 = $ws.Range("Q52").Value2
$R52 = $ws.Range("R52").Value2

$A53 = $ws.Range("A53").Value2
$B53 = $ws.Range("B53").Value2
$D53 = $ws.Range("D53").Value2
$E53 = $ws.Range("E53").Value2
$F53 = $ws.Range("F53").Value2
$G53 = $ws.Range("G53").Value2
$H53 = $ws.Range("H53").Value2
$P53 = $ws.Range("P53").Value2
$Q53 = $ws.Range("Q53").Value2
$R53 = $ws.Range("R53").Value2

$A54 = $ws.Range("A54").Value2
$B54 = $ws.Range("B54").Value2
$D54 = $ws.Range("D54").Value2
$E54 = $ws.Range("E54").Value2
$F54 = $ws.Range("F54").Value2
$G54 = $ws.Range("G54").Value2
$H54 = $ws.Range("H54").Value2
$P54 = $ws.Range("P54").Value2
$Q54 = $ws.Range("Q54").Value2
$R54 = $ws.Range("R54").Value2

$A55 = $ws.Range("A55").Value2
$B55 = $ws.Range("B55").Value2
$D55 = $ws.Range("D55").Value2
$E55 = $ws.Range("E55").Value2
$F55 = $ws.Range("F55").Value2
$G55 = $ws.Range("G55").Value2
$H55 = $ws.Range("H55").Value2
$P55 = $ws.Range("P55").Value2
$Q55 = $ws.Range("Q55").Value2
$R55 = $ws.Range("R55").Value2

$A56 = $ws.Range("A56").Value2
$B56 = $ws.Range("B56").Value2
$D56 = $ws.Range("D56").Value2
$E56 = $ws.Range("E56").Value2
$F56 = $ws.Range("F56").Value2
$G56 = $ws.Range("G56").Value2
$H56 = $ws.Range("H56").Value2
$P56 = $ws.Range("P56").Value2
$Q56 = $ws.Range("Q56").Value2
$R56 = $ws.Range("R56").Value2

$A57 = $ws.Range("A57").Value2
$B57 = $ws.Range("B57").Value2
$D57 = $ws.Range("D57").Value2
$E57 = $ws.Range("E57").Value2
$F57 = $ws.Range("F57").Value2
$G57 = $ws.Range("G57").Value2
$H57 = $ws.Range("H57").Value2
$P57 = $ws.Range("P57").Value2
$Q57 = $ws.Range("Q57").Value2
$R57 = $ws.Range("R57").Value2

$A59 = $ws.Range("A59").Value2
$B59 = $ws.Range("B59").Value2
$D59 = $ws.Range("D59").Value2
$E59 = $ws.Range("E59").Value2
$F59 = $ws.Range("F59").Value2
$G59 = $ws.Range("G59").Value2
$H59 = $ws.Range("H59").Value2
$P59 = $ws.Range("P59").Value2
$Q59 = $ws.Range("Q59").Value2
$R59 = $ws.Range("R59").Value2

$A60 = $ws.Range("A60").Value2
$B60 = $ws.Range("B60").Value2
$D60 = $ws.Range("D60").Value2
$E60 = $ws.Range("E60").Value2
$F60 = $ws.Range("F60").Value2
$G60 = $ws.Range("G60").Value2
$H60 = $ws.Range("H60").Value2
$P60 = $ws.Range("P60").Value2
$Q60 = $ws.Range("Q60").Value2
$R60 = $ws.Range("R60").Value2

# --- Write the corrected record into each row. ---
$ws.Range("A52").Value2 = $A54
$ws.Range("B52").Value2 = $B54
$ws.Range("D52").Value2 = $D54
$ws.Range("E52").Value2 = $E54
$ws.Range("F52").Value2 = $F54
$ws.Range("G52").Value2 = $G54
$ws.Range("H52").Value2 = $H54
$ws.Range("P52").Value2 = $P54
$ws.Range("Q52").Value2 = $Q54
$ws.Range("R52").Value2 = $R54

$ws.Range("A53").Value2 = $A57
$ws.Range("B53").Value2 = $B57
$ws.Range("D53").Value2 = $D57
$ws.Range("E53").Value2 = $E57
$ws.Range("F53").Value2 = $F57
$ws.Range("G53").Value2 = $G57
$ws.Range("H53").Value2 = $H57
$ws.Range("P53").Value2 = $P57
$ws.Range("Q53").Value2 = $Q57
$ws.Range("R53").Value2 = $R57

$ws.Range("A54").Value2 = $A52
$ws.Range("B54").Value2 = $B52
$ws.Range("D54").Value2 = $D52
$ws.Range("E54").Value2 = $E52
$ws.Range("F54").Value2 = $F52
$ws.Range("G54").Value2 = $G52
$ws.Range("H54").Value2 = $H52
$ws.Range("P54").Value2 = $P52
$ws.Range("Q54").Value2 = $Q52
$ws.Range("R54").Value2 = $R52

$ws.Range("A55").Value2 = $A53
$ws.Range("B55").Value2 = $B53
$ws.Range("D55").Value2 = $D53
$ws.Range("E55").Value2 = $E53
$ws.Range("F55").Value2 = $F53
$ws.Range("G55").Value2 = $G53
$ws.Range("H55").Value2 = $H53
$ws.Range("P55").Value2 = $P53
$ws.Range("Q55").Value2 = $Q53
$ws.Range("R55").Value2 = $R53

$ws.Range("A56").Value2 = $A55
$ws.Range("B56").Value2 = $B55
$ws.Range("D56").Value2 = $D55
$ws.Range("E56").Value2 = $E55
$ws.Range("F56").Value2 = $F55
$ws.Range("G56").Value2 = $G55
$ws.Range("H56").Value2 = $H55
$ws.Range("P56").Value2 = $P55
$ws.Range("Q56").Value2 = $Q55
$ws.Range("R56").Value2 = $R55

$ws.Range("A57").Value2 = $A56
$ws.Range("B57").Value2 = $B56
$ws.Range("D57").Value2 = $D56
$ws.Range("E57").Value2 = $E56
$ws.Range("F57").Value2 = $F56
$ws.Range("G57").Value2 = $G56
$ws.Range("H57").Value2 = $H56
$ws.Range("P57").Value2 = $P56
$ws.Range("Q57").Value2 = $Q56
$ws.Range("R57").Value2 = $R56

$ws.Range("A59").Value2 = $A60
$ws.Range("B59").Value2 = $B60
$ws.Range("D59").Value2 = $D60
$ws.Range("E59").Value2 = $E60
$ws.Range("F59").Value2 = $F60
$ws.Range("G59").Value2 = $G60
$ws.Range("H59").Value2 = $H60
$ws.Range("P59").Value2 = $P60
$ws.Range("Q59").Value2 = $Q60
$ws.Range("R59").Value2 = $R60

$ws.Range("A60").Value2 = $A59
$ws.Range("B60").Value2 = $B59
$ws.Range("D60").Value2 = $D59
$ws.Range("E60").Value2 = $E59
$ws.Range("F60").Value2 = $F59
$ws.Range("G60").Value2 = $G59
$ws.Range("H60").Value2 = $H59
$ws.Range("P60").Value2 = $P59
$ws.Range("Q60").Value2 = $Q59
$ws.Range("R60").Value2 = $R59

